$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weather-log rows to append (columns A-J):
# temp, temp_feels_like, pressure, humidity, description, clouds, place, wind, cloudiness, timestamp
$data = @(
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:34:22 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:34:34 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:34:46 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:34:57 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:35:08 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:35:19 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:35:30 02-12-2025"),
  @(12.47, 11.86, 1018, 80, "shower rain", 40, "Lisbon", 18.504, 40, "19:38:19 02-12-2025"),
  @(12.09, 11.52, 1019, 83, "few clouds",  20, "Lisbon", 18.504, 20, "19:52:22 02-12-2025"),
  @(12.09, 11.52, 1019, 83, "few clouds",  20, "Lisbon", 18.504, 20, "19:54:58 02-12-2025")
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $startRow + $i
    $rowValues = $data[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $colIndex = $c + 1
        $ws.Cells.Item($rowIndex, $colIndex).Value = $rowValues[$c]
    }
}
